$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Re-shuffle F:V blocks within same-kickoff-time groups ---
$groups = @(
    @(3,4),
    @(6,7,8,9),
    @(20,21),
    @(34,35,36),
    @(39,40,41,42),
    @(47,48,49,50),
    @(52,54),
    @(58,59,60),
    @(62,63,64,65,66)
)

$mapping = @{
    3 = 4;
    4 = 3;
    6 = 9;
    7 = 6;
    8 = 7;
    9 = 8;
    20 = 21;
    21 = 20;
    34 = 35;
    35 = 36;
    36 = 34;
    39 = 42;
    40 = 39;
    41 = 40;
    42 = 41;
    47 = 50;
    48 = 49;
    49 = 48;
    50 = 47;
    52 = 54;
    54 = 52;
    58 = 60;
    59 = 58;
    60 = 59;
    62 = 65;
    63 = 66;
    64 = 62;
    65 = 63;
    66 = 64;
}

foreach ($g in $groups) {
    $saved = @{}
    foreach ($r in $g) {
        $rowvals = @()
        for ($c = 6; $c -le 22; $c++) {
            $rowvals += ,$ws.Cells.Item($r, $c).Value2
        }
        $saved[$r] = $rowvals
    }
    foreach ($r in $g) {
        $src = $mapping[$r]
        $srcvals = $saved[$src]
        for ($c = 6; $c -le 22; $c++) {
            $ws.Cells.Item($r, $c).Value2 = $srcvals[$c-6]
        }
    }
}
# --- Step 2: Append new rows 68-72 ---
# Row 68
$ws.Range("A2").Copy()
$ws.Range("A68").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E68").PasteSpecial(-4122)
$ws.Range("A68").Value2 = 67
$ws.Range("B68").Value2 = 'greece'
$ws.Range("C68").Value2 = 'super-league-2'
$ws.Range("D68").Value2 = '2023-2024'
$ws.Range("E68").Value2 = 45241.54166666666
$ws.Range("F68").Value2 = 'Olympiacos Piraeus B'
$ws.Range("G68").Value2 = 1
$ws.Range("H68").Value2 = 'Diagoras'
$ws.Range("I68").Value2 = 0
$ws.Range("J68").Value2 = 1.68
$ws.Range("K68").Value2 = '10/11/2023 01:13'
$ws.Range("L68").Value2 = 1.22
$ws.Range("M68").Value2 = '11/11/2023 12:59'
$ws.Range("N68").Value2 = 3.44
$ws.Range("O68").Value2 = '10/11/2023 01:13'
$ws.Range("P68").Value2 = 5.93
$ws.Range("Q68").Value2 = '11/11/2023 12:59'
$ws.Range("R68").Value2 = 4.38
$ws.Range("S68").Value2 = '10/11/2023 01:13'
$ws.Range("T68").Value2 = 11.59
$ws.Range("U68").Value2 = '11/11/2023 12:59'
$ws.Range("V68").Value2 = 'https://www.betexplorer.com/football/greece/super-league-2/olympiacos-piraeus-diagoras-fc/ALEkYkkF/'

# Row 69
$ws.Range("A2").Copy()
$ws.Range("A69").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E69").PasteSpecial(-4122)
$ws.Range("A69").Value2 = 68
$ws.Range("B69").Value2 = 'greece'
$ws.Range("C69").Value2 = 'super-league-2'
$ws.Range("D69").Value2 = '2023-2024'
$ws.Range("E69").Value2 = 45241.58333333334
$ws.Range("F69").Value2 = 'Apollon Pontou'
$ws.Range("G69").Value2 = 2
$ws.Range("H69").Value2 = 'AEK Athens FC B'
$ws.Range("I69").Value2 = 2
$ws.Range("J69").Value2 = 2.52
$ws.Range("K69").Value2 = '10/11/2023 02:13'
$ws.Range("L69").Value2 = 3.36
$ws.Range("M69").Value2 = '11/11/2023 13:49'
$ws.Range("N69").Value2 = 2.82
$ws.Range("O69").Value2 = '10/11/2023 02:13'
$ws.Range("P69").Value2 = 2.73
$ws.Range("Q69").Value2 = '11/11/2023 13:49'
$ws.Range("R69").Value2 = 2.83
$ws.Range("S69").Value2 = '10/11/2023 02:13'
$ws.Range("T69").Value2 = 2.45
$ws.Range("U69").Value2 = '11/11/2023 13:49'
$ws.Range("V69").Value2 = 'https://www.betexplorer.com/football/greece/super-league-2/apollon-pontou-aek/APStLKRa/'

# Row 70
$ws.Range("A2").Copy()
$ws.Range("A70").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E70").PasteSpecial(-4122)
$ws.Range("A70").Value2 = 69
$ws.Range("B70").Value2 = 'greece'
$ws.Range("C70").Value2 = 'super-league-2'
$ws.Range("D70").Value2 = '2023-2024'
$ws.Range("E70").Value2 = 45241.58333333334
$ws.Range("F70").Value2 = 'Kalamata'
$ws.Range("G70").Value2 = 3
$ws.Range("H70").Value2 = 'Tilikratis L.'
$ws.Range("I70").Value2 = 0
$ws.Range("J70").Value2 = 1.13
$ws.Range("K70").Value2 = '11/11/2023 12:01'
$ws.Range("L70").Value2 = 1.13
$ws.Range("M70").Value2 = '11/11/2023 12:01'
$ws.Range("N70").Value2 = 7.63
$ws.Range("O70").Value2 = '11/11/2023 12:01'
$ws.Range("P70").Value2 = 7.63
$ws.Range("Q70").Value2 = '11/11/2023 12:01'
$ws.Range("R70").Value2 = 18.11
$ws.Range("S70").Value2 = '11/11/2023 12:01'
$ws.Range("T70").Value2 = 18.11
$ws.Range("U70").Value2 = '11/11/2023 12:01'
$ws.Range("V70").Value2 = 'https://www.betexplorer.com/football/greece/super-league-2/kalamata-tilikratis-lefkada/6m8bWTJR/'

# Row 71
$ws.Range("A2").Copy()
$ws.Range("A71").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E71").PasteSpecial(-4122)
$ws.Range("A71").Value2 = 70
$ws.Range("B71").Value2 = 'greece'
$ws.Range("C71").Value2 = 'super-league-2'
$ws.Range("D71").Value2 = '2023-2024'
$ws.Range("E71").Value2 = 45241.58333333334
$ws.Range("F71").Value2 = 'AEL Larissa'
$ws.Range("G71").Value2 = 2
$ws.Range("H71").Value2 = 'PAOK B'
$ws.Range("I71").Value2 = 0
$ws.Range("J71").Value2 = 1.43
$ws.Range("K71").Value2 = '10/11/2023 02:13'
$ws.Range("L71").Value2 = 1.33
$ws.Range("M71").Value2 = '11/11/2023 13:40'
$ws.Range("N71").Value2 = 3.99
$ws.Range("O71").Value2 = '10/11/2023 02:13'
$ws.Range("P71").Value2 = 4.69
$ws.Range("Q71").Value2 = '11/11/2023 13:40'
$ws.Range("R71").Value2 = 6.04
$ws.Range("S71").Value2 = '10/11/2023 02:13'
$ws.Range("T71").Value2 = 9.68
$ws.Range("U71").Value2 = '11/11/2023 13:40'
$ws.Range("V71").Value2 = 'https://www.betexplorer.com/football/greece/super-league-2/ael-larissa-paok/WYTxM0Cg/'

# Row 72
$ws.Range("A2").Copy()
$ws.Range("A72").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E72").PasteSpecial(-4122)
$ws.Range("A72").Value2 = 71
$ws.Range("B72").Value2 = 'greece'
$ws.Range("C72").Value2 = 'super-league-2'
$ws.Range("D72").Value2 = '2023-2024'
$ws.Range("E72").Value2 = 45241.58333333334
$ws.Range("F72").Value2 = 'Kampaniakos'
$ws.Range("G72").Value2 = 0
$ws.Range("H72").Value2 = 'Kozani FC'
$ws.Range("I72").Value2 = 1
$ws.Range("J72").Value2 = 2.61
$ws.Range("K72").Value2 = '10/11/2023 02:13'
$ws.Range("L72").Value2 = 2.94
$ws.Range("M72").Value2 = '11/11/2023 13:03'
$ws.Range("N72").Value2 = 2.87
$ws.Range("O72").Value2 = '10/11/2023 02:13'
$ws.Range("P72").Value2 = 2.79
$ws.Range("Q72").Value2 = '11/11/2023 13:03'
$ws.Range("R72").Value2 = 2.61
$ws.Range("S72").Value2 = '10/11/2023 02:13'
$ws.Range("T72").Value2 = 2.68
$ws.Range("U72").Value2 = '11/11/2023 13:03'
$ws.Range("V72").Value2 = 'https://www.betexplorer.com/football/greece/super-league-2/kampaniakos-kozani-fc/xCUYMtcm/'

$excel.CutCopyMode = $false